$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$c = $t.Cell(1, 1).Range
$c.Find.Execute("269×7=", $true, $false, $false, $false, $false, $true, 0, $false, "164×7=", 1) | Out-Null
$c = $t.Cell(1, 2).Range
$c.Find.Execute("558×4=", $true, $false, $false, $false, $false, $true, 0, $false, "718×5=", 1) | Out-Null
$c = $t.Cell(1, 3).Range
$c.Find.Execute("919×8=", $true, $false, $false, $false, $false, $true, 0, $false, "969×7=", 1) | Out-Null
$c = $t.Cell(1, 4).Range
$c.Find.Execute("941×8=", $true, $false, $false, $false, $false, $true, 0, $false, "927×4=", 1) | Out-Null
$c = $t.Cell(1, 5).Range
$c.Find.Execute("728×4=", $true, $false, $false, $false, $false, $true, 0, $false, "493×4=", 1) | Out-Null
$c = $t.Cell(5, 1).Range
$c.Find.Execute("949×8=", $true, $false, $false, $false, $false, $true, 0, $false, "579×9=", 1) | Out-Null
$c = $t.Cell(5, 2).Range
$c.Find.Execute("419×9=", $true, $false, $false, $false, $false, $true, 0, $false, "936×2=", 1) | Out-Null
$c = $t.Cell(5, 3).Range
$c.Find.Execute("897×8=", $true, $false, $false, $false, $false, $true, 0, $false, "518×5=", 1) | Out-Null
$c = $t.Cell(5, 4).Range
$c.Find.Execute("587×8=", $true, $false, $false, $false, $false, $true, 0, $false, "929×9=", 1) | Out-Null
$c = $t.Cell(5, 5).Range
$c.Find.Execute("114×8=", $true, $false, $false, $false, $false, $true, 0, $false, "303×6=", 1) | Out-Null
$c = $t.Cell(10, 1).Range
$c.Find.Execute("563×6=", $true, $false, $false, $false, $false, $true, 0, $false, "757×4=", 1) | Out-Null
$c = $t.Cell(10, 2).Range
$c.Find.Execute("672×4=", $true, $false, $false, $false, $false, $true, 0, $false, "384×3=", 1) | Out-Null
$c = $t.Cell(10, 3).Range
$c.Find.Execute("391×5=", $true, $false, $false, $false, $false, $true, 0, $false, "210×6=", 1) | Out-Null
$c = $t.Cell(10, 4).Range
$c.Find.Execute("328×5=", $true, $false, $false, $false, $false, $true, 0, $false, "971×5=", 1) | Out-Null
$c = $t.Cell(10, 5).Range
$c.Find.Execute("523×8=", $true, $false, $false, $false, $false, $true, 0, $false, "220×8=", 1) | Out-Null
$c = $t.Cell(15, 1).Range
$c.Find.Execute("977×8=", $true, $false, $false, $false, $false, $true, 0, $false, "808×2=", 1) | Out-Null
$c = $t.Cell(15, 2).Range
$c.Find.Execute("352×6=", $true, $false, $false, $false, $false, $true, 0, $false, "747×7=", 1) | Out-Null
$c = $t.Cell(15, 3).Range
$c.Find.Execute("232×5=", $true, $false, $false, $false, $false, $true, 0, $false, "588×6=", 1) | Out-Null
$c = $t.Cell(15, 4).Range
$c.Find.Execute("359×5=", $true, $false, $false, $false, $false, $true, 0, $false, "470×3=", 1) | Out-Null
$c = $t.Cell(15, 5).Range
$c.Find.Execute("600×4=", $true, $false, $false, $false, $false, $true, 0, $false, "651×8=", 1) | Out-Null
$c = $t.Cell(20, 1).Range
$c.Find.Execute("219×9=", $true, $false, $false, $false, $false, $true, 0, $false, "776×2=", 1) | Out-Null
$c = $t.Cell(20, 2).Range
$c.Find.Execute("735×9=", $true, $false, $false, $false, $false, $true, 0, $false, "305×2=", 1) | Out-Null
$c = $t.Cell(20, 3).Range
$c.Find.Execute("918×3=", $true, $false, $false, $false, $false, $true, 0, $false, "262×8=", 1) | Out-Null
$c = $t.Cell(20, 4).Range
$c.Find.Execute("114×8=", $true, $false, $false, $false, $false, $true, 0, $false, "983×6=", 1) | Out-Null
$c = $t.Cell(20, 5).Range
$c.Find.Execute("296×5=", $true, $false, $false, $false, $false, $true, 0, $false, "189×3=", 1) | Out-Null
